# Update milestone dates / status across Showroom_Progress, Warehouse_Progress
# and Settings sheets to reflect the new Dec-2025 timeline (see commit message).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell while keeping date-look-alike text cells
# stored as TEXT (not auto-converted to Excel date serials). Non date-like
# strings / numbers are written normally.
# ---------------------------------------------------------------------------
function Set-TextValue {
    param(
        $Range,
        [string]$Value
    )
    $Range.NumberFormat = "@"
    $Range.Value = $Value
}

# ---------------------------------------------------------------------------
# Sheet: Showroom_Progress
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Showroom_Progress")

Set-TextValue $ws.Range("B2") "2025-03-01"
Set-TextValue $ws.Range("E2") "2025-02-25"

Set-TextValue $ws.Range("B3") "2025-05-15"
Set-TextValue $ws.Range("E3") "2025-05-10"

Set-TextValue $ws.Range("B4") "2025-07-15"
$ws.Range("C4").Value = 50
$ws.Range("F4").Value = "Paused - reviewing options"

Set-TextValue $ws.Range("B5") "2026-03-01"
$ws.Range("F5").Value = "On hold pending strategic review"

Set-TextValue $ws.Range("B6") "2026-05-15"
$ws.Range("F6").Value = "Custom displays to be ordered"

Set-TextValue $ws.Range("B7") "2026-07-01"
$ws.Range("F7").Value = "Revised target date"

# ---------------------------------------------------------------------------
# Sheet: Warehouse_Progress
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Warehouse_Progress")

Set-TextValue $ws.Range("B2") "2025-06-01"
Set-TextValue $ws.Range("E2") "2025-05-28"

Set-TextValue $ws.Range("B3") "2025-09-01"
Set-TextValue $ws.Range("E3") "2025-08-25"

Set-TextValue $ws.Range("B4") "2026-04-01"
$ws.Range("F4").Value = "Started 4 Nov 2025 - landlord managing refurbishment (4 weeks in)"

Set-TextValue $ws.Range("B5") "2026-05-01"

Set-TextValue $ws.Range("B6") "2026-06-15"

Set-TextValue $ws.Range("B7") "2026-07-15"

# ---------------------------------------------------------------------------
# Sheet: Settings
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Settings")

Set-TextValue $ws.Range("B3") "2026-07-01"
Set-TextValue $ws.Range("B6") "2026-07-15"
